# Rotate the "Recorded By" (column G) comma-separated list for every data
# row so the last entry moves to the front, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "a, b, c" -> "c, a, b"
# Single-value cells (no comma) and empty cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $val.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $count = $trimmed.Count
    $last = $trimmed[$count - 1]
    $rest = $trimmed[0..($count - 2)]
    $newOrder = @($last) + $rest

    $newVal = [string]::Join(", ", $newOrder)
    $cell.Value = $newVal
}
